$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02, 1.02)
$CValues = @(1.028382726690068, 1.030325156207344, 1.031575916713346, 1.032100302485043, 1.032188265985439, 1.031582929173871, 1.029040466392303, 1.024512010649539, 1.021458538040969, 1.020127711578314, 1.01963204376406, 1.019738427457443, 1.020086767081912, 1.020301211910651, 1.021546674941731, 1.022325576679331, 1.022779064742774, 1.022933552496067, 1.022242094251393, 1.019984227046583, 1.018556845553122, 1.019314277503414, 1.022279818936753, 1.025688655906901)
$DValues = @(1.031956381641022, 1.033444566559109, 1.034402243388289, 1.03480360995001, 1.034870929061682, 1.034407611312612, 1.032460427482943, 1.02898777636902, 1.02664332777225, 1.025620862726974, 1.025239947503323, 1.025321706544415, 1.025589399298964, 1.025754183534658, 1.026711029079981, 1.027309257450311, 1.027657491255012, 1.027776111663344, 1.027245146097292, 1.025510601778107, 1.024413494438553, 1.024995720621448, 1.02727411742821, 1.029890603266122)
$EValues = @(1.028104918743613, 1.029806281502341, 1.030901306094183, 1.031360276552743, 1.031437259709995, 1.030907444256973, 1.028681134605332, 1.024711875118401, 1.022032857138948, 1.020864622950147, 1.020429421781899, 1.02052283180804, 1.020828675135185, 1.021016946365465, 1.022110213241758, 1.022793770502555, 1.023191688920832, 1.023327236037993, 1.022720513133791, 1.020738647155386, 1.019485221330257, 1.020150394360884, 1.022753617421858, 1.025743681221349)
$FValues = @(1.026932373238568, 1.029524038627695, 1.031193934884516, 1.031894310293911, 1.032011811068504, 1.031203299751472, 1.027809734143545, 1.021773441640437, 1.017708309824393, 1.015937672331922, 1.015278355770071, 1.015419855815476, 1.015883206473769, 1.016168475258108, 1.017825596880881, 1.018862237401202, 1.019465892055272, 1.019671554811429, 1.018751119566556, 1.015746806536324, 1.013848455577009, 1.014855720324384, 1.018801332020944, 1.023340958592643)
$IValues = @(1.031040219857644, 1.031561218825675, 1.031894252084095, 1.032033290655612, 1.032056579333374, 1.03189611371675, 1.031217146423547, 1.029988953776024, 1.029148153623308, 1.028778711876489, 1.028640665187868, 1.028670313967416, 1.028767317688402, 1.028826975892427, 1.029172558031801, 1.029387885970067, 1.029512965892232, 1.029555527568042, 1.029364836910874, 1.028738775238941, 1.028340397064731, 1.028552039314348, 1.029375253386271, 1.030310301287518)
$JValues = @(1.033535352532168, 1.035113283085584, 1.036128069497138, 1.036553218528711, 1.036624517657347, 1.036133756088438, 1.034069929831923, 1.030384198110756, 1.027892378560218, 1.026804783745611, 1.026399473740668, 1.026486474794945, 1.02677130800013, 1.026946626030518, 1.027964374157586, 1.028600449105458, 1.028970631418185, 1.029096714343712, 1.028532290318682, 1.026687468677798, 1.025519850498669, 1.026139568821822, 1.028563090902105, 1.031343030326792)
$KValues = @(1.034763002514901, 1.036058040677981, 1.036890255375326, 1.037238760779641, 1.037297197382472, 1.036894917427127, 1.035201872853842, 1.032173399892303, 1.030122693147082, 1.02922687481471, 1.028892921344701, 1.028964610547932, 1.029199294892947, 1.029343730938973, 1.030181977962747, 1.030705667613679, 1.031010371515885, 1.031114140364525, 1.030649559003193, 1.029130219735514, 1.028167955501564, 1.028678742366845, 1.030674914402184, 1.032961820654317)
$LValues = @(1.03092270516103, 1.032429523201192, 1.033398233619777, 1.033803996914007, 1.033872040170093, 1.033403661239287, 1.03143326196043, 1.027911721328044, 1.025529112329011, 1.024488757943796, 1.024100988773694, 1.024184227532332, 1.024456732253611, 1.024624453482798, 1.025597971791251, 1.026206290206441, 1.026560277864123, 1.026680837999662, 1.026141109916578, 1.024376523530688, 1.023259316095616, 1.023852313861153, 1.026170564660129, 1.028828153701717)
$MValues = @(1.029753576994986, 1.032148041176021, 1.033690114697445, 1.034336697554263, 1.034445162956073, 1.033698761009374, 1.030564336566534, 1.024983137572472, 1.021220561254592, 1.019580779263869, 1.01897004916031, 1.019101128046595, 1.019530329944396, 1.019794556136921, 1.021329160968484, 1.022288911694161, 1.022847702800694, 1.023038065580387, 1.022186044933655, 1.019403986461453, 1.017645267233159, 1.018578519526368, 1.02223252912944, 1.026433164878064)

$startRow = 2
for ($i = 0; $i -lt 24; $i++) {
    $row = $startRow + $i
    $ws.Range("B" + $row).Value = $BValues[$i]
    $ws.Range("C" + $row).Value = $CValues[$i]
    $ws.Range("D" + $row).Value = $DValues[$i]
    $ws.Range("E" + $row).Value = $EValues[$i]
    $ws.Range("F" + $row).Value = $FValues[$i]
    $ws.Range("I" + $row).Value = $IValues[$i]
    $ws.Range("J" + $row).Value = $JValues[$i]
    $ws.Range("K" + $row).Value = $KValues[$i]
    $ws.Range("L" + $row).Value = $LValues[$i]
    $ws.Range("M" + $row).Value = $MValues[$i]
}
